$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3804.5
$ws.Range("I74").Value = 3698
$ws.Range("J74").Value = 3831.125
$ws.Range("K74").Value = 3698
$ws.Range("L74").Value = 3831.125
$ws.Range("M74").Value = -2762
$ws.Range("N74").Value = -5703.125
$ws.Range("H77").Value = 3804.5
$ws.Range("I77").Value = 3698
$ws.Range("J77").Value = 3831.125
$ws.Range("K77").Value = 18490
$ws.Range("L77").Value = 19155.625
$ws.Range("M77").Value = -13810
$ws.Range("N77").Value = -28515.625
$ws.Range("H140").Value = 52800
$ws.Range("J140").Value = 52800
$ws.Range("L140").Value = 52800
$ws.Range("N140").Value = -63160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7776.3
$ws.Range("I32").Value = 8907.333000000001
$ws.Range("J32").Value = 3252.1667
$ws.Range("K32").Value = 8907.333000000001
$ws.Range("L32").Value = 3252.1667
$ws.Range("M32").Value = -8620.333000000001
$ws.Range("N32").Value = -3826.1667
$ws.Range("H88").Value = 2567.0667
$ws.Range("I88").Value = 2764.8333
$ws.Range("J88").Value = 2435.2222
$ws.Range("K88").Value = 2764.8333
$ws.Range("L88").Value = 2435.2222
$ws.Range("M88").Value = -2358.8333
$ws.Range("N88").Value = -3247.2222
$ws.Range("H91").Value = 2567.0667
$ws.Range("I91").Value = 2764.8333
$ws.Range("J91").Value = 2435.2222
$ws.Range("K91").Value = 2764.8333
$ws.Range("L91").Value = 2435.2222
$ws.Range("M91").Value = -1360.8333
$ws.Range("N91").Value = -5243.2222
$ws.Range("H132").Value = 4238904.5
$ws.Range("I132").Value = 5001333
$ws.Range("K132").Value = 15003999
$ws.Range("M132").Value = -15001469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19232498
$ws.Range("I86").Value = 1722.4
$ws.Range("J86").Value = 83335080
$ws.Range("K86").Value = 1722.4
$ws.Range("L86").Value = 83335080
$ws.Range("M86").Value = -599.4000000000001
$ws.Range("N86").Value = -83337326
$ws.Range("H89").Value = 19232498
$ws.Range("I89").Value = 1722.4
$ws.Range("J89").Value = 83335080
$ws.Range("K89").Value = 8612
$ws.Range("L89").Value = 416675400
$ws.Range("M89").Value = -2996
$ws.Range("N89").Value = -416686632
$ws.Range("H105").Value = 4048.093
$ws.Range("I105").Value = 3118.0625
$ws.Range("J105").Value = 4599.222
$ws.Range("K105").Value = 3118.0625
$ws.Range("L105").Value = 4599.222
$ws.Range("M105").Value = -1371.0625
$ws.Range("N105").Value = -8093.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2402.0557
$ws.Range("I58").Value = 885.375
$ws.Range("J58").Value = 3615.4
$ws.Range("K58").Value = 885.375
$ws.Range("L58").Value = 3615.4
$ws.Range("M58").Value = -682.375
$ws.Range("N58").Value = -4021.4
$ws.Range("H107").Value = 415.66666
$ws.Range("I107").Value = 603.75
$ws.Range("J107").Value = 164.88889
$ws.Range("K107").Value = 603.75
$ws.Range("L107").Value = 164.88889
$ws.Range("M107").Value = 1316.25
$ws.Range("N107").Value = -4004.88889
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180
$ws.Range("H132").Value = 6967.3184
$ws.Range("I132").Value = 7765.25
$ws.Range("J132").Value = 4839.5
$ws.Range("K132").Value = 23295.75
$ws.Range("L132").Value = 14518.5
$ws.Range("M132").Value = -20765.75
$ws.Range("N132").Value = -19578.5
$ws.Range("H136").Value = 2402.0557
$ws.Range("I136").Value = 885.375
$ws.Range("J136").Value = 3615.4
$ws.Range("K136").Value = 2656.125
$ws.Range("L136").Value = 10846.2
$ws.Range("M136").Value = -106.125
$ws.Range("N136").Value = -15946.2
$ws.Range("H140").Value = 29560
$ws.Range("J140").Value = 29560
$ws.Range("L140").Value = 29560
$ws.Range("N140").Value = -39920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1469
$ws.Range("I5").Value = 955.2857
$ws.Range("K5").Value = 2865.8571
$ws.Range("M5").Value = -2753.8571
$ws.Range("H131").Value = 821.97
$ws.Range("I131").Value = 539.8
$ws.Range("J131").Value = 836.82104
$ws.Range("K131").Value = 1619.4
$ws.Range("L131").Value = 2510.46312
$ws.Range("M131").Value = 3420.6
$ws.Range("N131").Value = -12590.46312
$ws.Range("H135").Value = 1469
$ws.Range("I135").Value = 955.2857
$ws.Range("K135").Value = 8597.5713
$ws.Range("M135").Value = -6062.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20657.46
$ws.Range("I70").Value = 27393.54
$ws.Range("J70").Value = 7185.3076
$ws.Range("K70").Value = 27393.54
$ws.Range("L70").Value = 7185.3076
$ws.Range("M70").Value = -27123.54
$ws.Range("N70").Value = -7725.3076
$ws.Range("H73").Value = 20657.46
$ws.Range("I73").Value = 27393.54
$ws.Range("J73").Value = 7185.3076
$ws.Range("K73").Value = 27393.54
$ws.Range("L73").Value = 7185.3076
$ws.Range("M73").Value = -26457.54
$ws.Range("N73").Value = -9057.3076
$ws.Range("H80").Value = 16777696
$ws.Range("J80").Value = 2553927.2
$ws.Range("L80").Value = 2553927.2
$ws.Range("N80").Value = -2555923.2
$ws.Range("H83").Value = 16777696
$ws.Range("J83").Value = 2553927.2
$ws.Range("L83").Value = 12769636
$ws.Range("N83").Value = -12779620
$ws.Range("H138").Value = 58056.57
$ws.Range("J138").Value = 58056.57
$ws.Range("L138").Value = 58056.57
$ws.Range("N138").Value = -68336.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1781.6818
$ws.Range("J82").Value = 1834.5454
$ws.Range("L82").Value = 1834.5454
$ws.Range("N82").Value = -2556.5454
$ws.Range("H85").Value = 1781.6818
$ws.Range("J85").Value = 1834.5454
$ws.Range("L85").Value = 1834.5454
$ws.Range("N85").Value = -4330.5454
$ws.Range("H93").Value = 1550.7142
$ws.Range("I93").Value = 1387.1428
$ws.Range("J93").Value = 1714.2858
$ws.Range("K93").Value = 1387.1428
$ws.Range("L93").Value = 1714.2858
$ws.Range("M93").Value = -139.1428000000001
$ws.Range("N93").Value = -4210.2858
$ws.Range("H136").Value = 10748.786
$ws.Range("I136").Value = 2407.0908
$ws.Range("K136").Value = 7221.2724
$ws.Range("M136").Value = -4671.2724
$ws.Range("H139").Value = 45217.637
$ws.Range("J139").Value = 45674.4
$ws.Range("L139").Value = 45674.4
$ws.Range("N139").Value = -55954.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 60465
$ws.Range("J125").Value = 60465
$ws.Range("L125").Value = 60465
$ws.Range("N125").Value = -70305
$ws.Range("H136").Value = 1570.3
$ws.Range("I136").Value = 1691.8334
$ws.Range("J136").Value = 1388
$ws.Range("K136").Value = 5075.5002
$ws.Range("L136").Value = 4164
$ws.Range("M136").Value = -2525.5002
$ws.Range("N136").Value = -9264
$ws.Range("H138").Value = 51949.668
$ws.Range("I138").Value = 25000
$ws.Range("K138").Value = 25000
$ws.Range("M138").Value = -19860

Write-Output "done"